$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that will receive new text values - reset to Normal style + Text format
# first so Excel does not auto-convert numeric-looking strings (e.g. "1.010")
# into numbers, then reset the style back to Normal (no explicit style) afterwards
# so the saved cell format matches the original (unstyled) cells.

$targetCells = @("D2", "E2", "D3", "E3", "D4", "E4", "B5", "C5", "D5", "E5", "B6", "C6", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "B38", "C38", "D38", "E38", "B39", "C39", "D39", "E39", "B40", "C40", "D40", "E40", "B41", "C41", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "E45", "B46", "C46", "D46", "E46", "B47", "C47", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).Style = "Normal"
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new values
$ws.Range('D2').Value = '27.081.75'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '1.794.77'
$ws.Range('E3').Value = '  -2.23%  '
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('B5').Value = 'USDC'
$ws.Range('C5').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D5').Value = '1.008'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '308.18'
$ws.Range('E6').Value = '  -1.86%  '
$ws.Range('D7').Value = '0.4185'
$ws.Range('E7').Value = '  -1.21%  '
$ws.Range('D8').Value = '0.3562'
$ws.Range('E8').Value = '  -3.19%  '
$ws.Range('D9').Value = '0.07069'
$ws.Range('E9').Value = '  -2.57%  '
$ws.Range('D10').Value = '0.8436'
$ws.Range('E10').Value = '  -2.43%  '
$ws.Range('D11').Value = '20.04'
$ws.Range('E11').Value = '  -3.00%  '
$ws.Range('D12').Value = '1.848.87'
$ws.Range('E12').Value = '  -5.17%  '
$ws.Range('D13').Value = '5.268'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('D14').Value = '6.332'
$ws.Range('E14').Value = '  -2.61%  '
$ws.Range('D15').Value = '0.06851'
$ws.Range('E15').Value = '  -1.84%  '
$ws.Range('D16').Value = '1.011'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').Value = '79.93'
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('E18').Value = '  -3.01%  '
$ws.Range('D19').Value = '1.008'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').Value = '15.04'
$ws.Range('E20').Value = '  -1.65%  '
$ws.Range('D21').Value = '27.186.91'
$ws.Range('E21').Value = '  -3.00%  '
$ws.Range('D22').Value = '5.039'
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').Value = '10.69'
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('D24').Value = '2.047.02'
$ws.Range('E24').Value = '  -4.06%  '
$ws.Range('D25').Value = '1.960'
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('D26').Value = '152.87'
$ws.Range('E26').Value = '  -0.79%  '
$ws.Range('D27').Value = '18.16'
$ws.Range('E27').Value = '  -1.50%  '
$ws.Range('D28').Value = '5.001'
$ws.Range('E28').Value = '  -4.40%  '
$ws.Range('D29').Value = '112.71'
$ws.Range('E29').Value = '  -2.21%  '
$ws.Range('D30').Value = '1.662'
$ws.Range('E30').Value = '  -8.94%  '
$ws.Range('D31').Value = '0.08887'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').Value = '0.7232'
$ws.Range('E32').Value = '  -5.68%  '
$ws.Range('D33').Value = '2.878'
$ws.Range('E33').Value = '  -2.80%  '
$ws.Range('D34').Value = '4.353'
$ws.Range('E34').Value = '  -3.95%  '
$ws.Range('D35').Value = '1.009'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').Value = '1.075'
$ws.Range('E36').Value = '  -5.39%  '
$ws.Range('D37').Value = '1.070'
$ws.Range('E37').Value = '  -3.13%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01889'
$ws.Range('E38').Value = '  -2.53%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05098'
$ws.Range('E39').Value = '  -4.71%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.4941'
$ws.Range('E40').Value = '  -2.53%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '0.1615'
$ws.Range('E41').Value = '  -2.69%  '
$ws.Range('D42').Value = '2.672'
$ws.Range('E42').Value = '  -5.48%  '
$ws.Range('D43').Value = '6.148'
$ws.Range('E43').Value = '  -9.28%  '
$ws.Range('D44').Value = '8.007'
$ws.Range('E44').Value = '  -5.10%  '
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '10.22'
$ws.Range('E46').Value = '  -2.70%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '104.51'
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('D48').Value = '0.06317'
$ws.Range('E48').Value = '  -3.38%  '
$ws.Range('D49').Value = '0.4527'
$ws.Range('E49').Value = '  -2.92%  '
$ws.Range('D50').Value = '1.588'
$ws.Range('E50').Value = '  -1.60%  '
$ws.Range('D51').Value = '62.08'
$ws.Range('E51').Value = '  -3.88%  '

# Restore default (unstyled) appearance now that values are text-safe
foreach ($addr in $targetCells) {
    $ws.Range($addr).Style = "Normal"
}
